$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Teams sheet: add a new team "Esteban Cordero" and sort the team list
# ---------------------------------------------------------------------------
$teams = $wb.Worksheets.Item("Teams")
$teams.Range("A6").Value = "Esteban Cordero"

$sortObj = $teams.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($teams.Range("A2:A6")) | Out-Null
$sortObj.SetRange($teams.Range("A2:B6"))
$sortObj.Header = 0
$sortObj.Apply()

$teams.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 2) Matches sheet: the first round's fixtures were reshuffled / reset
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches")

# Match 1 (row 2): away team now FGJ, kickoff still TBD, score reset
$matches.Range("F2").Value = "FGJ"
$matches.Range("A2").NumberFormat = "h:mm"
$matches.Range("A2").Value = "Por definir"
$matches.Range("B2").Value = "Por definir"
$matches.Range("G2").Value = 0
$matches.Range("H2").Value = 0

# Match 2 (row 3): away team now the new team, Esteban Cordero, score reset
$matches.Range("F3").Value = "Esteban Cordero"
$matches.Range("A3").NumberFormat = "h:mm"
$matches.Range("A3").Value = "Por definir"
$matches.Range("B3").NumberFormat = "h:mm"
$matches.Range("B3").Value = "Por definir"
$matches.Range("G3").Value = 0
$matches.Range("H3").Value = 0

# Match 3 (row 4): stadium/home team now La 40, score + video cleared
$matches.Range("D4").Value = "La 40"
$matches.Range("E4").Value = "La 40"
$matches.Range("A4").NumberFormat = "h:mm"
$matches.Range("A4").Value = "Por definir"
$matches.Range("B4").NumberFormat = "h:mm"
$matches.Range("B4").Value = "Por definir"
$matches.Range("G4:H4").ClearContents()
$matches.Range("I4").Value = ""

# Match 4 (row 5): fixture removed entirely, only the time-formatted
# placeholders remain
$matches.Range("A5").NumberFormat = "h:mm"
$matches.Range("A5:F5").ClearContents()

$matches.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 3) Videos sheet: update the embed link shown for the "Rolitas" video
# ---------------------------------------------------------------------------
$videos = $wb.Worksheets.Item("Videos")
$videos.Range("B4").Value = "https://www.youtube.com/embed/VQZr2vG5k_I"
$videos.Range("B7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Scorers sheet: only the remembered selection moved
# ---------------------------------------------------------------------------
$scorers = $wb.Worksheets.Item("Scorers")
$scorers.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) Restore the remembered selections on Teams/Matches and make Teams the
#    active sheet (it was Matches before the edit)
# ---------------------------------------------------------------------------
$matches.Range("E7").Select() | Out-Null

$teams.Activate()
$teams.Range("D10").Select() | Out-Null
